$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 5, column G (was 4.000.000, should match 2.000.000 like column F)
$ws.Range("G5").Value = "2.000.000"

# Add new contractor rows (7, 8, 9)
$ws.Range("A7").Value = "NT004"
$ws.Range("B7").Value = "Nhà thầu 4"
$ws.Range("C7").Value = "01-01-2025"
$ws.Range("D7").Value = "8.000.000"
$ws.Range("E7").Value = "01-01-2027"
$ws.Range("F7").Value = "3.000.000"
$ws.Range("G7").Value = "4.000.000"
$ws.Range("H7").Value = "400.000"

$ws.Range("A8").Value = "NT005"
$ws.Range("B8").Value = "Nhà thầu 5"
$ws.Range("C8").Value = "01-01-2025"
$ws.Range("D8").Value = "9.000.000"
$ws.Range("E8").Value = "01-01-2027"
$ws.Range("F8").Value = "3.000.000"
$ws.Range("G8").Value = "5.000.000"
$ws.Range("H8").Value = "400.000"

$ws.Range("A9").Value = "NT006"
$ws.Range("B9").Value = "Nhà thầu 6"
$ws.Range("C9").Value = "01-01-2025"
$ws.Range("D9").Value = "10.000.000"
$ws.Range("E9").Value = "01-01-2027"
$ws.Range("F9").Value = "3.000.000"
$ws.Range("G9").Value = "6.000.000"
$ws.Range("H9").Value = "400.000"

# Update the active selection to match the new edit location
$ws.Range("G12").Select()

$wb.Save()
